# Update the "Population Mean Sigma Unknown" worksheet's sample inputs and
# fix the "Choosing your Sample Size" callout heading, then restore that
# sheet as the active sheet/selection (matching the author's final view
# state of zoom 130%, cell H3 selected).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Population Mean Sigma Unknown")

# Fix the header text in E3 (was "Choosing your theSample Size").
$ws1.Range("E3").Value = "Choosing your Sample Size"

# Update the sample statistics that drive the confidence-interval formulas.
$ws1.Range("C6").Value = 200      # Sample average
$ws1.Range("C7").Value = 56       # Sample standard deviation (s)
$ws1.Range("C8").Value = 18       # Sample size (n)

# Make this sheet active again, with the zoom/selection the author left it at.
$ws1.Activate()
$ws1.Range("H3").Select()
$excel.ActiveWindow.Zoom = 130
